$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "№"
$ws.Range("B1").Value = "Наименование"

# Data rows - new content replacing old
$ws.Range("A2").Value = 4
$ws.Range("B2").Value = "* Манжета М60х80 ГОСТ 22704"

$ws.Range("A3").Value = 14
$ws.Range("B3").Value = "* Рукав газосварочный I-6,3-0,63-У ГОСТ 9356"

$ws.Range("A4").Value = 19
$ws.Range("B4").Value = "* Гвоздь строительный круглый головка плоская 3х80 ГОСТ 4028"

# The old rows 5-8 no longer hold data; clear them out entirely
$ws.Range("A5:B8").ClearContents()

# Update selection to match the new active cell
$ws.Range("B11").Select()
